$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "Y5LR9J"
$ws.Range("B63").Value = "Film de fusor Brother"
$ws.Range("C63").Value = "HL L5000 L5100 L5200, DCP L5500 5580 5585 5590 L5600 L5650 8150 8155, MFC L5700 L5750 L5755 L5800 L5850 L5900 L6200 L6250 L6300 L6400 L6700 L6750 L6800 L6900 8530 8535 8540"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 150000
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Formula = "=(E63-D63)*G63"
$ws.Range("I63").Formula = "=D63*F63"
$ws.Range("J63").Value = 0
